# "wrapping up test file audit"
#
# 1. optimization_parameters: remove the stray "Sheet" row (row 16, values
#    3 / 4) that had been left in the sheet; everything below shifts up by
#    one row. Removing the row also drops the now-unused "Sheet" shared
#    string (and its associated integer-format style) automatically.
# 2. threshold_b becomes the active / selected sheet (with A2 selected),
#    and the previously-selected optimization_parameters sheet keeps a
#    selection on the row that slid up into the deleted row's place.

$wb = $excel.ActiveWorkbook

$wsParams = $wb.Worksheets.Item("optimization_parameters")
$wsParams.Rows.Item(16).Delete()
$wsParams.Rows.Item(16).Select()

$wsThreshold = $wb.Worksheets.Item("threshold_b")
$wsThreshold.Activate()
$wsThreshold.Range("A2").Select()
